# إضافة حدث جديد في Card11
#
# The "Card11" service-log sheet keeps one row per maintenance event.
# Row 16's placeholder cells (columns B:K and P) were left truly blank;
# this backfills them with the literal text "nan" used elsewhere in the
# sheet for not-applicable values, and appends a brand-new event as row 17
# (date 14\8\2024 — machine greased completely + serviced, by "تيم العمل").

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Card11")

# Reference an existing, unstyled data cell so newly-written cells pick up
# the same ("no special formatting") style as the rest of the table.
$plainStyle = $ws.Range("A16").Style

# --- Backfill row 16's previously-empty cells with "nan" ---
$ws.Range("B16:K16").Value = "nan"
$ws.Range("P16").Value = "nan"

# --- Append the new event as row 17 ---
# Card number (text "11", matching the rest of the column).
$ws.Range("A17").Value = "'11"
$ws.Range("A17").Style = $plainStyle

# Empty/non-applicable cells for this event (kept as blank text cells,
# matching the sheet's convention rather than truly-empty/undefined cells).
$ws.Range("B17:K17").Value = "'"
$ws.Range("B17:K17").Style = $plainStyle
$ws.Range("M17").Value = "'"
$ws.Range("M17").Style = $plainStyle
$ws.Range("P17").Value = "'"
$ws.Range("P17").Style = $plainStyle

# The actual new event data.
$ws.Range("L17").Value = "14\8\2024"
$ws.Range("N17").Value = "تم تشحيم المكنه بالكامل +عمل صيانه"
$ws.Range("O17").Value = "تيم العمل"
